$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 797.4286
$ws.Range("I4").Value = 736
$ws.Range("J4").Value = 951
$ws.Range("K4").Value = 736
$ws.Range("L4").Value = 951
$ws.Range("M4").Value = -622
$ws.Range("N4").Value = -1179
# Row 11
$ws.Range("H11").Value = 27.947369
$ws.Range("I11").Value = 27.947369
$ws.Range("K11").Value = 27.947369
$ws.Range("M11").Value = 112.052631
# Row 17
$ws.Range("H17").Value = 2047.3928
$ws.Range("J17").Value = 2187
$ws.Range("L17").Value = 6561
$ws.Range("N17").Value = -6897
# Row 28
$ws.Range("H28").Value = 10488.4
$ws.Range("I28").Value = 925.3333
$ws.Range("J28").Value = 24833
$ws.Range("K28").Value = 925.3333
$ws.Range("L28").Value = 24833
$ws.Range("M28").Value = -440.3333
$ws.Range("N28").Value = -25803
# Row 39
$ws.Range("H39").Value = 171.5
$ws.Range("I39").Value = 45.545456
$ws.Range("K39").Value = 136.636368
$ws.Range("M39").Value = 159.363632
# Row 116
$ws.Range("H116").Value = 3999.5
$ws.Range("I116").Value = 3999.5
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 3999.5
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -557.5
$ws.Range("N116").ClearContents()
# Row 125
$ws.Range("H125").Value = 2439
$ws.Range("I125").Value = 1798.75
$ws.Range("K125").Value = 16188.75
$ws.Range("M125").Value = -13728.75
# Row 132
$ws.Range("H132").Value = 12393.88
$ws.Range("I132").Value = 11367.4
$ws.Range("K132").Value = 34102.2
$ws.Range("M132").Value = -31572.2
# Row 137
$ws.Range("H137").Value = 2624.3
$ws.Range("I137").Value = 1570.5
$ws.Range("J137").Value = 3828.6428
$ws.Range("K137").Value = 4711.5
$ws.Range("L137").Value = 11485.9284
$ws.Range("M137").Value = -2161.5
$ws.Range("N137").Value = -16585.9284

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 56
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("N56").ClearContents()
# Row 61
$ws.Range("H61").Value = 4174.6665
$ws.Range("I61").Value = 3749.625
$ws.Range("J61").Value = 5024.75
$ws.Range("K61").Value = 3749.625
$ws.Range("L61").Value = 5024.75
$ws.Range("M61").Value = -3537.625
$ws.Range("N61").Value = -5448.75
# Row 122
$ws.Range("H122").Value = 1343.2727
$ws.Range("I122").Value = 1388.8
$ws.Range("J122").Value = 888
$ws.Range("K122").Value = 4166.4
$ws.Range("L122").Value = 2664
$ws.Range("M122").Value = -1716.4
$ws.Range("N122").Value = -7564
# Row 132
$ws.Range("H132").Value = 2649.5833
$ws.Range("I132").Value = 1310.6666
$ws.Range("K132").Value = 3931.9998
$ws.Range("M132").Value = -1401.9998
# Row 136
$ws.Range("H136").Value = 4174.6665
$ws.Range("I136").Value = 3749.625
$ws.Range("J136").Value = 5024.75
$ws.Range("K136").Value = 11248.875
$ws.Range("L136").Value = 15074.25
$ws.Range("M136").Value = -8698.875
$ws.Range("N136").Value = -20174.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 16
$ws.Range("H16").Value = 1200
$ws.Range("I16").Value = 1200
$ws.Range("K16").Value = 1200
$ws.Range("M16").Value = -1030
# Row 94
$ws.Range("H94").Value = 899
$ws.Range("I94").Value = 899
$ws.Range("K94").Value = 899
$ws.Range("M94").Value = -448
# Row 134
$ws.Range("H134").Value = 3336.652
$ws.Range("I134").Value = 2591.7
$ws.Range("K134").Value = 7775.099999999999
$ws.Range("M134").Value = -5240.099999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 52
$ws.Range("H52").Value = 303910
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
# Row 58
$ws.Range("H58").Value = 4028.611
$ws.Range("I58").Value = 3116.6924
$ws.Range("K58").Value = 3116.6924
$ws.Range("M58").Value = -2913.6924
# Row 94
$ws.Range("H94").Value = 3828.9167
$ws.Range("I94").Value = 1776.2858
$ws.Range("K94").Value = 1776.2858
$ws.Range("M94").Value = -1325.2858
# Row 125
$ws.Range("H125").Value = 86666.336
$ws.Range("J125").Value = 86666.336
$ws.Range("L125").Value = 86666.336
$ws.Range("N125").Value = -91586.336
# Row 136
$ws.Range("H136").Value = 4028.611
$ws.Range("I136").Value = 3116.6924
$ws.Range("K136").Value = 9350.0772
$ws.Range("M136").Value = -6800.0772

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 75
$ws.Range("H75").Value = 2249.6667
$ws.Range("J75").Value = 2249.6667
$ws.Range("L75").Value = 6749.000100000001
$ws.Range("N75").Value = -8745.000100000001
# Row 78
$ws.Range("H78").Value = 2249.6667
$ws.Range("J78").Value = 2249.6667
$ws.Range("L78").Value = 20247.0003
$ws.Range("N78").Value = -30231.0003
# Row 108
$ws.Range("H108").Value = 1403.8334
$ws.Range("I108").Value = 918.5
$ws.Range("K108").Value = 2755.5
$ws.Range("M108").Value = 124.5
# Row 122
$ws.Range("H122").Value = 1117.4
$ws.Range("I122").Value = 998
$ws.Range("J122").Value = 1197
$ws.Range("K122").Value = 8982
$ws.Range("L122").Value = 10773
$ws.Range("M122").Value = -6532
$ws.Range("N122").Value = -15673
# Row 129
$ws.Range("H129").Value = 2034.0834
$ws.Range("I129").Value = 696.8333
$ws.Range("J129").Value = 3371.3333
$ws.Range("K129").Value = 2090.4999
$ws.Range("L129").Value = 10113.9999
$ws.Range("M129").Value = 2909.5001
$ws.Range("N129").Value = -20113.9999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 20000
$ws.Range("J43").Value = 20000
$ws.Range("L43").Value = 20000
$ws.Range("N43").Value = -20302
# Row 46
$ws.Range("H46").Value = 27999.666
$ws.Range("J46").Value = 31999.5
$ws.Range("L46").Value = 31999.5
$ws.Range("N46").Value = -32311.5
# Row 94
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
# Row 122
$ws.Range("H122").Value = 297610.28
$ws.Range("I122").Value = 360295.56
$ws.Range("K122").Value = 1080886.68
$ws.Range("M122").Value = -1078436.68
# Row 132
$ws.Range("H132").Value = 37511.695
$ws.Range("I132").Value = 47769
$ws.Range("J132").Value = 5457.625
$ws.Range("K132").Value = 143307
$ws.Range("L132").Value = 16372.875
$ws.Range("M132").Value = -140777
$ws.Range("N132").Value = -21432.875

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1461.3846
$ws.Range("J22").Value = 1777.7778
$ws.Range("L22").Value = 1777.7778
$ws.Range("N22").Value = -2367.7778
# Row 27
$ws.Range("H27").Value = 1461.3846
$ws.Range("J27").Value = 1777.7778
$ws.Range("L27").Value = 1777.7778
$ws.Range("N27").Value = -1991.7778
# Row 93
$ws.Range("H93").Value = 2132.9
$ws.Range("I93").Value = 2132.9
$ws.Range("K93").Value = 2132.9
$ws.Range("M93").Value = -884.9000000000001
# Row 132
$ws.Range("H132").Value = 5193.1577
$ws.Range("I132").Value = 4627.353
$ws.Range("J132").Value = 10002.5
$ws.Range("K132").Value = 13882.059
$ws.Range("L132").Value = 30007.5
$ws.Range("M132").Value = -11352.059
$ws.Range("N132").Value = -35067.5
# Row 141
$ws.Range("H141").Value = 77236
$ws.Range("J141").Value = 76695
$ws.Range("L141").Value = 76695
$ws.Range("N141").Value = -87055

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 3375
$ws.Range("I2").Value = 4333.3335
$ws.Range("K2").Value = 4333.3335
$ws.Range("M2").Value = -4221.3335
# Row 4
$ws.Range("H4").Value = 16377.857
$ws.Range("I4").Value = 34675
$ws.Range("J4").Value = 2655
$ws.Range("K4").Value = 34675
$ws.Range("L4").Value = 2655
$ws.Range("M4").Value = -34562
$ws.Range("N4").Value = -2881
# Row 5
$ws.Range("H5").Value = 15909932
$ws.Range("I5").Value = 20714536
$ws.Range("J5").Value = 7501875
$ws.Range("K5").Value = 20714536
$ws.Range("L5").Value = 7501875
$ws.Range("M5").Value = -20714424
$ws.Range("N5").Value = -7502099
# Row 39
$ws.Range("H39").Value = 60024.75
$ws.Range("J39").Value = 60024.75
$ws.Range("L39").Value = 60024.75
$ws.Range("N39").Value = -60850.75
# Row 82
$ws.Range("H82").Value = 35500
$ws.Range("J82").Value = 35500
$ws.Range("L82").Value = 35500
$ws.Range("N82").Value = -36266
# Row 85
$ws.Range("H85").Value = 35500
$ws.Range("J85").Value = 35500
$ws.Range("L85").Value = 35500
$ws.Range("N85").Value = -38152
# Row 132
$ws.Range("H132").Value = 2702.4138
$ws.Range("I132").Value = 2513.25
$ws.Range("K132").Value = 7539.75
$ws.Range("M132").Value = -5009.75
